$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1904761904761905
$ws.Range("C2").Value = 0.5625
$ws.Range("J2").Value = 0.0119047619047619
$ws.Range("O2").Value = 0.002976190476190476
$ws.Range("P2").Value = 0.1517857142857143
$ws.Range("S2").Value = 0.08035714285714286
$ws.Range("B3").Value = 0.005181347150259068
$ws.Range("C3").Value = 0.0155440414507772
$ws.Range("J3").Value = 0.02590673575129534
$ws.Range("P3").Value = 0.7564766839378239
$ws.Range("S3").Value = 0.1968911917098446
$ws.Range("J4").Value = 0.02173913043478261
$ws.Range("P4").Value = 0.6956521739130435
$ws.Range("S4").Value = 0.2826086956521739
$ws.Range("B6").Value = 0.06276150627615062
$ws.Range("D6").Value = 0.01255230125523013
$ws.Range("F6").Value = 0.07949790794979079
$ws.Range("J6").Value = 0.3096234309623431
$ws.Range("O6").Value = 0.01673640167364017
$ws.Range("Q6").Value = 0.08368200836820083
$ws.Range("R6").Value = 0.04184100418410042
$ws.Range("S6").Value = 0.393305439330544
$ws.Range("B7").Value = 0.1084905660377359
$ws.Range("D7").Value = 0.01886792452830189
$ws.Range("F7").Value = 0.0330188679245283
$ws.Range("J7").Value = 0.1179245283018868
$ws.Range("O7").Value = 0.02830188679245283
$ws.Range("Q7").Value = 0.169811320754717
$ws.Range("R7").Value = 0.1037735849056604
$ws.Range("S7").Value = 0.419811320754717
$ws.Range("B8").Value = 0.1021276595744681
$ws.Range("D8").Value = 0.01276595744680851
$ws.Range("F8").Value = 0.05106382978723404
$ws.Range("J8").Value = 0.151063829787234
$ws.Range("O8").Value = 0.01276595744680851
$ws.Range("Q8").Value = 0.1574468085106383
$ws.Range("R8").Value = 0.07872340425531915
$ws.Range("S8").Value = 0.4340425531914894
$ws.Range("B9").Value = 0.1683168316831683
$ws.Range("D9").Value = 0.01485148514851485
$ws.Range("F9").Value = 0.05445544554455446
$ws.Range("J9").Value = 0.1386138613861386
$ws.Range("O9").Value = 0.0297029702970297
$ws.Range("Q9").Value = 0.08415841584158416
$ws.Range("R9").Value = 0.103960396039604
$ws.Range("S9").Value = 0.405940594059406
$ws.Range("B10").Value = 0.1012658227848101
$ws.Range("D10").Value = 0.02131912058627582
$ws.Range("E10").Value = 0.001332445036642239
$ws.Range("F10").Value = 0.05929380413057962
$ws.Range("J10").Value = 0.128580946035976
$ws.Range("O10").Value = 0.01932045303131246
$ws.Range("Q10").Value = 0.2111925383077948
$ws.Range("R10").Value = 0.08927381745502998
$ws.Range("S10").Value = 0.3684210526315789
$ws.Range("G11").Value = 0.1518624641833811
$ws.Range("J11").Value = 0.1002865329512894
$ws.Range("K11").Value = 0.2063037249283668
$ws.Range("L11").Value = 0.5214899713467048
$ws.Range("S11").Value = 0.02005730659025788
$ws.Range("G12").Value = 0.6878306878306878
$ws.Range("J12").Value = 0.2592592592592592
$ws.Range("K12").Value = 0.01058201058201058
$ws.Range("L12").Value = 0.02116402116402116
$ws.Range("S12").Value = 0.02116402116402116
$ws.Range("G13").Value = 0.6595744680851063
$ws.Range("J13").Value = 0.2340425531914894
$ws.Range("S13").Value = 0.1063829787234043
$ws.Range("F15").Value = 0.03162055335968379
$ws.Range("H15").Value = 0.1541501976284585
$ws.Range("I15").Value = 0.05928853754940711
$ws.Range("J15").Value = 0.3794466403162055
$ws.Range("K15").Value = 0.05928853754940711
$ws.Range("M15").Value = 0.01185770750988142
$ws.Range("O15").Value = 0.06719367588932806
$ws.Range("S15").Value = 0.2371541501976284
$ws.Range("F16").Value = 0.0179372197309417
$ws.Range("H16").Value = 0.1479820627802691
$ws.Range("I16").Value = 0.07174887892376682
$ws.Range("J16").Value = 0.4798206278026906
$ws.Range("K16").Value = 0.08520179372197309
$ws.Range("M16").Value = 0.0179372197309417
$ws.Range("O16").Value = 0.07623318385650224
$ws.Range("S16").Value = 0.1031390134529148
$ws.Range("F17").Value = 0.03260869565217391
$ws.Range("H17").Value = 0.1413043478260869
$ws.Range("I17").Value = 0.1304347826086956
$ws.Range("J17").Value = 0.4217391304347826
$ws.Range("K17").Value = 0.1065217391304348
$ws.Range("M17").Value = 0.01304347826086956
$ws.Range("N17").Value = 0.002173913043478261
$ws.Range("O17").Value = 0.05434782608695652
$ws.Range("S17").Value = 0.09782608695652174
$ws.Range("F18").Value = 0.01809954751131222
$ws.Range("H18").Value = 0.1764705882352941
$ws.Range("I18").Value = 0.05882352941176471
$ws.Range("J18").Value = 0.4298642533936652
$ws.Range("K18").Value = 0.1131221719457014
$ws.Range("M18").Value = 0.01357466063348416
$ws.Range("O18").Value = 0.06334841628959276
$ws.Range("S18").Value = 0.1266968325791855
$ws.Range("F19").Value = 0.02217453505007153
$ws.Range("H19").Value = 0.2110157367668097
$ws.Range("I19").Value = 0.07010014306151645
$ws.Range("J19").Value = 0.3819742489270386
$ws.Range("K19").Value = 0.1144492131616595
$ws.Range("M19").Value = 0.02217453505007153
$ws.Range("N19").Value = 0.00357653791130186
$ws.Range("O19").Value = 0.07010014306151645
$ws.Range("S19").Value = 0.1044349070100143
